# Generate Report for Handoff
# Refresh the localization-status report: Priority values move from
# "low" to "ht", and the handoff timestamps for the zh-cn / de-de
# languages (and the rolled-up "Latest HO Xliff Generate Date" on the
# Overview sheet) are updated to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-20 16:41:54"
$wsOverview.Range("G5").Value = "2016-08-20 16:41:54"
$wsOverview.Range("G6").Value = "2016-08-20 16:41:54"
$wsOverview.Range("G7").Value = "2016-08-20 16:41:54"

# --- zh-cn sheet: Priority (E) and Latest Handoff Datetime (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("H4").Value = "2016-08-20 16:41:49"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("H5").Value = "2016-08-20 16:41:49"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("H6").Value = "2016-08-20 16:41:49"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("H7").Value = "2016-08-20 16:41:49"

# --- de-de sheet: Priority (E) and Latest Handoff Datetime (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("H4").Value = "2016-08-20 16:41:54"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("H5").Value = "2016-08-20 16:41:54"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("H6").Value = "2016-08-20 16:41:54"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("H7").Value = "2016-08-20 16:41:54"
